# career_analysis.xlsx :: "Vibematch Questions" sheet
#
# Adds a dedicated "Options" column (new D) describing the Likert scale /
# response type for each question, and pushes the existing "Algorithmic
# Impact" commentary into a new column E with refreshed wording that calls
# out the 1-5 scoring weight explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vibematch Questions")

# --- Column widths -------------------------------------------------------
# Raw OOXML <col width> ends up ~0.8333 (5/6) wider than the COM
# ColumnWidth value for this workbook's default font/theme, so back that
# offset out to land exactly on the target raw width of 60 for both D & E
# (matches the existing column B which is also width 60).
$targetRawWidth = 60
$comColumnWidth = $targetRawWidth - (5 / 6)
$ws.Columns.Item(4).ColumnWidth = $comColumnWidth
$ws.Columns.Item(5).ColumnWidth = $comColumnWidth

# --- Header row ------------------------------------------------------------
$ws.Range("D1").Value = "Options"
$ws.Range("E1").Value = "Algorithmic Impact"

# --- Likert rows (v_01..v_14): rows 2-15 -----------------------------------
# Each row's old column-D "Maps to {...}. Contributes to 40% RIASEC Score."
# text moves to column E with the tail reworded; column D gets the shared
# 1-5 Likert legend.
$likertOptions = "Strongly Disagree (1), Disagree (2), Neutral (3), Agree (4), Strongly Agree (5)"
$oldSuffix = "Contributes to 40% RIASEC Score."
$newSuffix = "Score (1-5) weighted by 40% RIASEC component."

for ($row = 2; $row -le 15; $row++) {
    $oldImpact = $ws.Range("D$row").Text
    $newImpact = $oldImpact.Replace($oldSuffix, $newSuffix)

    $ws.Range("D$row").Value = $likertOptions
    $ws.Range("E$row").Value = $newImpact
}

# --- Final free-text row (v_15): row 16 -------------------------------------
$ws.Range("D16").Value = "Text Response"
$ws.Range("E16").Value = "Scoring (Practical): Bonus if text matches career/bucket (Positive Reinforcement)."
